$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# These "Price" cells are stored as plain text in the workbook (note the
# "." used as a thousands separator elsewhere in the column). Several of the
# new values look like ordinary decimals, so without forcing a text format
# Excel would silently convert them to numbers. Mark them as Text first.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply the updated values from the diff
$ws.Range("D2").Value = "37.706.39"
$ws.Range("E2").Value = "  +0.98%  "
$ws.Range("D3").Value = "2.075.88"
$ws.Range("E3").Value = "  +0.51%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "232.53"
$ws.Range("E5").Value = "  -0.47%  "
$ws.Range("E6").Value = "  +0.64%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").Value = "56.97"
$ws.Range("E8").Value = "  +0.27%  "
$ws.Range("E9").Value = "  +1.10%  "
$ws.Range("D10").Value = "0.0786"
$ws.Range("E10").Value = "  +3.01%  "
$ws.Range("E11").Value = "  +2.97%  "
$ws.Range("D12").Value = "2.368.13"
$ws.Range("E12").Value = "  -0.07%  "
$ws.Range("D13").Value = "14.40"
$ws.Range("E13").Value = "  +0.06%  "
$ws.Range("D14").Value = "20.87"
$ws.Range("E14").Value = "  +0.81%  "
$ws.Range("D15").Value = "0.758"
$ws.Range("E15").Value = "  -2.04%  "
$ws.Range("D16").Value = "5.24"
$ws.Range("E16").Value = "  +1.66%  "
$ws.Range("D17").Value = "2.070.98"
$ws.Range("E17").Value = "  +0.22%  "
$ws.Range("D18").Value = "37.629.14"
$ws.Range("E18").Value = "  +0.94%  "
$ws.Range("D19").Value = "6.17"
$ws.Range("E19").Value = "  -3.43%  "
$ws.Range("D20").Value = "70.79"
$ws.Range("E20").Value = "  +1.85%  "
$ws.Range("D21").Value = "0.0₃0818"
$ws.Range("E21").Value = "  +0.65%  "
$ws.Range("D22").Value = "227.56"
$ws.Range("E22").Value = "  +0.75%  "
$ws.Range("E23").Value = "  +0.00%  "
$ws.Range("E24").Value = "  -0.93%  "
$ws.Range("E25").Value = "  -0.65%  "
$ws.Range("D26").Value = "169.70"
$ws.Range("E26").Value = "  +2.00%  "
$ws.Range("D27").Value = "0.139"
$ws.Range("E27").Value = "  +10.55%  "
$ws.Range("D28").Value = "8.88"
$ws.Range("E28").Value = "  +1.21%  "
$ws.Range("D29").Value = "1.43"
$ws.Range("E29").Value = "  -0.01%  "
$ws.Range("D30").Value = "19.34"
$ws.Range("E30").Value = "  +2.05%  "
$ws.Range("E31").Value = "  +0.51%  "
$ws.Range("D32").Value = "4.61"
$ws.Range("E32").Value = "  +2.84%  "
$ws.Range("D33").Value = "0.0623"
$ws.Range("E33").Value = "  +1.03%  "
$ws.Range("D34").Value = "4.59"
$ws.Range("E34").Value = "  -0.31%  "
$ws.Range("D35").Value = "2.49"
$ws.Range("E35").Value = "  +0.38%  "
$ws.Range("E36").Value = "  +3.32%  "
$ws.Range("D37").Value = "3.36"
$ws.Range("E37").Value = "  +4.51%  "
$ws.Range("E38").Value = "  +0.06%  "
$ws.Range("D39").Value = "5.40"
$ws.Range("E39").Value = "  -4.27%  "
$ws.Range("D40").Value = "0.0986"
$ws.Range("E40").Value = "  +5.78%  "
$ws.Range("B41").Value = "HuobiToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D41").Value = "2.94"
$ws.Range("E41").Value = "  -0.60%  "
$ws.Range("B42").Value = "Aave"
$ws.Range("C42").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D42").Value = "98.63"
$ws.Range("E42").Value = "  +2.76%  "
$ws.Range("D43").Value = "4.38"
$ws.Range("E43").Value = "  +3.47%  "
$ws.Range("E44").Value = "  +1.07%  "
$ws.Range("D45").Value = "1.453.78"
$ws.Range("E45").Value = "  -1.32%  "
$ws.Range("E46").Value = "  -1.26%  "
$ws.Range("D47").Value = "1.05"
$ws.Range("E47").Value = "  +2.70%  "
$ws.Range("D48").Value = "15.56"
$ws.Range("E48").Value = "  +2.18%  "
$ws.Range("D49").Value = "7.39"
$ws.Range("E49").Value = "  +3.07%  "
$ws.Range("E50").Value = "  +1.38%  "
$ws.Range("D51").Value = "47.34"
$ws.Range("E51").Value = "  +7.84%  "
